$d = $word.ActiveDocument

# The first table in the document is the project-info header table.
# Its first row gets an explicit height, and the first (label) cell of
# each of its four rows is given a solid blue (accent2, 50% darker) fill.
$t = $d.Tables(1)

# Row 1: set an explicit row height of 1275 twips (63.75 points).
$row1 = $t.Rows(1)
$row1.Height = 63.75

$fillColor = 8015618  # RGB(0x02, 0x4f, 0x7a) -> 024f7a

for ($i = 1; $i -le $t.Rows.Count; $i++) {
  $cell = $t.Cell($i, 1)
  $cell.Shading.Texture = "clear"
  $cell.Shading.ForegroundPatternColor = $fillColor
  $cell.Shading.BackgroundPatternColor = $fillColor
}
